$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# 1. "Logs" sheet: append row 3 with the second test mail entry
# ----------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A3").Value = "Kun je 10 dozen schroeven bestellen?"
$logs.Range("B3").Value = "mailmind.test@zohomail.eu"
$logs.Range("C3").Value = "Testmail #2: Kun je 10 dozen schroeven bestellen?"
$logs.Range("D3").Value = "Bestelling / Levering"
$logs.Range("E3").Value = "Geachte klant,`nBedankt voor uw e-mail. Helaas kan ik geen bestellingen plaatsen, maar ik kan u doorverwijzen naar het bestelteam binnen ons bedrijf. Graag ontvang ik de contactgegevens van uw bedrijf, zodat ik de juiste persoon met u in contact kan brengen.`nIk zie uw reactie graag tegemoet.`nMet vriendelijke groet,`n[Naam]`nE-mailassistent"
$logs.Range("F3").Value = "2025-06-29 13:52:05"
$logs.Range("G3").Value = "Ja"
$logs.Range("H3").Value = "Ja"
$logs.Range("I3").Value = "Nee"

# Extend the conditional-formatting blocks so they also cover row 3
$logs.Range("D2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D3"))
$logs.Range("G2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G3"))
$logs.Range("H2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H3"))
$logs.Range("I2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I3"))

# ----------------------------------------------------------------------
# 2. "Dashboard" sheet: append row 3 with the tally for the new category
# ----------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A3").Value = "Bestelling / Levering"
$dash.Range("B3").Value = 1

# ----------------------------------------------------------------------
# 3. Chart on the Dashboard sheet: extend the category/value series refs
# ----------------------------------------------------------------------
$chart = $dash.ChartObjects().Item(1).Chart
$ser = $chart.SeriesCollection().Item(1)
$ser.XValues = "'Dashboard'!`$A`$2:`$A`$3"
$ser.Values = "'Dashboard'!`$B`$2:`$B`$3"
